$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new rows after row 17 (the current last data row) to make room for
#    the two new worker records. This shifts the old rows 18-23 down to 20-25,
#    preserving the signature block merges (B22:C22 etc -> B24:C24 etc).
$ws.Range("18:19").Insert()

# 2. Copy row 17's formatting (font, fill, borders, number formats, alignment)
#    down into the two freshly inserted rows so the new records look the same
#    as the existing table rows.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Populate the two new worker rows.
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143392269"
$ws.Range("D18").Value = "LINA YULIANA MARIN LOPEZ"
$ws.Range("E18").Value = "2509"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1041328245"
$ws.Range("D19").Value = "MONICA YOHANA MARIN LOPEZ"
$ws.Range("E19").Value = "2509"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# 4. Update the summary fields at the top of the statement.
#    Valor Mora total = sum of the four worker rows.
$ws.Range("E11").Value = 124922
#    Cant. Trabajadores (distinct workers) and Cant. Periodos (distinct periods).
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 2

Write-Host "done"
